# Generate Report for Handback
#
# After a locale's translated files are handed back "in sync" with en-US,
# the per-locale status report grows two columns - "Latest Target File" (F)
# and "Latest Handback File" (G) - and the "Latest Handback DateTime" (H) /
# "Status" columns get populated with the real handback data instead of
# their handoff-time placeholders.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# 1. Status column: "Ready for handoff" -> "Handed back: in sync with en-US"
#    (Overview sheet B/C, and each locale sheet's Status column C)
# ---------------------------------------------------------------------
$ovw = $wb.Worksheets.Item("Overview")
$ovw.Range("B2").Value = $newStatus
$ovw.Range("C2").Value = $newStatus
$ovw.Range("B3").Value = $newStatus
$ovw.Range("C3").Value = $newStatus

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

# ---------------------------------------------------------------------
# 2. zh-cn sheet: add "Latest Target File" (F) / "Latest Handback File" (G)
#    hyperlinks, mirroring the source (A) / handoff (D) file for each row,
#    and fill in the real handback datetime (H).
# ---------------------------------------------------------------------
$zhcn.Hyperlinks.Add(
    $zhcn.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/d1dbb7eb75d96c139a1570fb114e9c6e4865cc75/e2e/af782635-8cd8-428e-b60e-7ba862439f1f.md",
    "",
    "",
    "af782635-8cd8-428e-b60e-7ba862439f1f.md"
) | Out-Null

$zhcn.Hyperlinks.Add(
    $zhcn.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/2365941ea3c6c1c36e85e077f8d01cfa83aebf71/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/hb/af782635-8cd8-428e-b60e-7ba862439f1f.2b5828a45aea85cdb0f06cf5eabed8d8b4149e76.zh-cn.xlf",
    "",
    "",
    "af782635-8cd8-428e-b60e-7ba862439f1f.2b5828a45aea85cdb0f06cf5eabed8d8b4149e76.zh-cn.xlf"
) | Out-Null

$zhcn.Hyperlinks.Add(
    $zhcn.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/d1dbb7eb75d96c139a1570fb114e9c6e4865cc75/e2e/c24ac456-c25d-494b-9f17-797937471c65.md",
    "",
    "",
    "c24ac456-c25d-494b-9f17-797937471c65.md"
) | Out-Null

$zhcn.Hyperlinks.Add(
    $zhcn.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/2365941ea3c6c1c36e85e077f8d01cfa83aebf71/ol-handback/OpenLocalizationTestOrg/oltest-zhcn-fly/yuwzho/hb/c24ac456-c25d-494b-9f17-797937471c65.e3306c3a115cc7c6fb777a212affef074f2a58d2.zh-cn.xlf",
    "",
    "",
    "c24ac456-c25d-494b-9f17-797937471c65.e3306c3a115cc7c6fb777a212affef074f2a58d2.zh-cn.xlf"
) | Out-Null

$zhcn.Range("H2").Value = "2016-03-20 17:06:59"
$zhcn.Range("H3").Value = "2016-03-20 17:06:59"

# ---------------------------------------------------------------------
# 3. de-de sheet: same shape of change as zh-cn, with its own handback
#    timestamp and de-de handoff filenames.
# ---------------------------------------------------------------------
$dede.Hyperlinks.Add(
    $dede.Range("F2"),
    "https://github.com/OpenLocalizationTest/oltest/blob/d1dbb7eb75d96c139a1570fb114e9c6e4865cc75/e2e/af782635-8cd8-428e-b60e-7ba862439f1f.md",
    "",
    "",
    "af782635-8cd8-428e-b60e-7ba862439f1f.md"
) | Out-Null

$dede.Hyperlinks.Add(
    $dede.Range("G2"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/623183d076ea3bfce62aba547db47487c821016f/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/hb/af782635-8cd8-428e-b60e-7ba862439f1f.2b5828a45aea85cdb0f06cf5eabed8d8b4149e76.de-de.xlf",
    "",
    "",
    "af782635-8cd8-428e-b60e-7ba862439f1f.2b5828a45aea85cdb0f06cf5eabed8d8b4149e76.de-de.xlf"
) | Out-Null

$dede.Hyperlinks.Add(
    $dede.Range("F3"),
    "https://github.com/OpenLocalizationTest/oltest/blob/d1dbb7eb75d96c139a1570fb114e9c6e4865cc75/e2e/c24ac456-c25d-494b-9f17-797937471c65.md",
    "",
    "",
    "c24ac456-c25d-494b-9f17-797937471c65.md"
) | Out-Null

$dede.Hyperlinks.Add(
    $dede.Range("G3"),
    "https://github.com/OpenLocalizationTestOrg/olhandback/blob/623183d076ea3bfce62aba547db47487c821016f/ol-handback/OpenLocalizationTestOrg/oltest-dede-fly/yuwzho/hb/c24ac456-c25d-494b-9f17-797937471c65.e3306c3a115cc7c6fb777a212affef074f2a58d2.de-de.xlf",
    "",
    "",
    "c24ac456-c25d-494b-9f17-797937471c65.e3306c3a115cc7c6fb777a212affef074f2a58d2.de-de.xlf"
) | Out-Null

$dede.Range("H2").Value = "2016-03-20 17:07:13"
$dede.Range("H3").Value = "2016-03-20 17:07:13"

Write-Output "Handback report generated."
